# Markowitz workbook - "Optimization v0 - functional - YA"
# Re-derives the I/J/K/L return series from rows 7:9 instead of 6:9
# (row 6 held a placeholder zero return that's now removed), switches the
# annualization factor from 250 to 252 trading days, and rebuilds the
# covariance matrix (Q5:T8) off the shorter I7:I9 .. L7:L9 ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the placeholder zero-return row (I6:L6) ------------------
# AVERAGE / VAR.S / COVAR formulas elsewhere already reference I6:I9 etc,
# so once I6:L6 are blank those formulas naturally recompute over I7:I9.
$ws.Range("I6:L6").ClearContents()

# --- 2. L7 / L8 : re-enter explicitly (was part of a J7:L7 / L8:L9 fill)
$ws.Range("L7").Formula = "=G7/G6-1"
$ws.Range("L8").Formula = "=G8/G7-1"

# --- 3. Rebuild the covariance matrix Q5:T8 off I7:I9 .. L7:L9 --------
$ws.Range("Q5").Formula = "=COVAR(`$I`$7:`$I`$9,I`$7:I`$9)"
$ws.Range("R5").Formula = "=COVAR(`$I`$7:`$I`$9,J`$7:J`$9)"
$ws.Range("S5").Formula = "=COVAR(`$I`$7:`$I`$9,K`$7:K`$9)*250"
$ws.Range("T5").Formula = "=COVAR(`$I`$7:`$I`$9,L`$7:L`$9)*250"

$ws.Range("Q6").Formula = "=COVAR(`$J`$7:`$J`$9,I`$7:I`$9)"
$ws.Range("R6:T6").Formula = "=COVAR(`$J`$7:`$J`$9,J`$7:J`$9)*250"

$ws.Range("Q7").Formula = "=COVAR(`$K`$7:`$K`$9,I`$7:I`$9)*250"
$ws.Range("R7:T7").Formula = "=COVAR(`$K`$7:`$K`$9,J`$7:J`$9)*250"

$ws.Range("Q8").Formula = "=COVAR(`$L`$7:`$L`$9,I`$7:I`$9)*250"
$ws.Range("R8:S8").Formula = "=COVAR(`$L`$7:`$L`$9,J`$7:J`$9)*250"
$ws.Range("T8").Formula = "=COVAR(`$L`$7:`$L`$9,L`$7:L`$9)*250"

# --- 4. Give the rebuilt covariance matrix a finer number format ------
$ws.Range("Q5:T8").NumberFormat = "0.0000000"
$ws.Range("Q5:T8").HorizontalAlignment = -4152  # xlRight

# Re-apply the bold highlight font to the same cells that carried it
# before (diagonal / "K" row emphasis), matching the author's layout.
$ws.Range("Q5:T5").Font.Bold = $false
$ws.Range("Q6").Font.Bold = $false
$ws.Range("R6:T6").Font.Bold = $true
$ws.Range("Q7:T7").Font.Bold = $true
$ws.Range("Q8:T8").Font.Bold = $false

# --- 5. Annualize off 252 trading days instead of 250 (I17:L17) -------
$ws.Range("I17").Formula = "=I14*252"
$ws.Range("J17:L17").Formula = "=J14*252"

# --- 6. Column widths (Q:T got wider / a new col) ----------------------
$ws.Columns("Q:R").ColumnWidth = 10.76
$ws.Columns("S").ColumnWidth = 10.42
$ws.Columns("T").ColumnWidth = 8.92

# --- 7. View state: scroll position + active selection -----------------
$ws.Select()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("R30").Select()

$excel.CalculateFullRebuild()
